# Adds 5 new fungal observation records (rows 8-12) to the worksheet,
# replicating the sparse column layout used by the existing data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: column letter, kind (num/str/bool), value.
# "str" values are prefixed with a leading apostrophe so Excel stores them
# as literal text (this also prevents date-looking strings such as
# "2023-10-03" from being auto-converted into date serial numbers).
$newRows = @(
    @{
        Row = 8
        Cells = @(
            @{ Col = "A"; Kind = "num"; Val = 112501198 }
            @{ Col = "B"; Kind = "num"; Val = 89834 }
            @{ Col = "C"; Kind = "str"; Val = '''Ovaliderad' }
            @{ Col = "D"; Kind = "str"; Val = '''NT' }
            @{ Col = "E"; Kind = "num"; Val = 658 }
            @{ Col = "F"; Kind = "str"; Val = '''Rosenticka' }
            @{ Col = "G"; Kind = "str"; Val = '''Rhodofomes roseus' }
            @{ Col = "H"; Kind = "str"; Val = '''(Alb. & Schwein.) Kotl. & Pouzar' }
            @{ Col = "I"; Kind = "str"; Val = '''' }
            @{ Col = "J"; Kind = "str"; Val = '''' }
            @{ Col = "K"; Kind = "str"; Val = '''' }
            @{ Col = "N"; Kind = "str"; Val = '''' }
            @{ Col = "P"; Kind = "str"; Val = '''Erkinjänkkä (Erkinjänkkä), T lm' }
            @{ Col = "Q"; Kind = "num"; Val = 749922 }
            @{ Col = "R"; Kind = "num"; Val = 7535992 }
            @{ Col = "S"; Kind = "num"; Val = 10 }
            @{ Col = "T"; Kind = "str"; Val = '''Norrbotten' }
            @{ Col = "U"; Kind = "str"; Val = '''Kiruna' }
            @{ Col = "V"; Kind = "str"; Val = '''Torne lappmark' }
            @{ Col = "W"; Kind = "str"; Val = '''Jukkasjärvi' }
            @{ Col = "Y"; Kind = "str"; Val = '''2023-10-03' }
            @{ Col = "AA"; Kind = "str"; Val = '''2023-10-03' }
            @{ Col = "AD"; Kind = "bool"; Val = $false }
            @{ Col = "AE"; Kind = "bool"; Val = $false }
            @{ Col = "AF"; Kind = "str"; Val = '''' }
            @{ Col = "AG"; Kind = "bool"; Val = $false }
            @{ Col = "AT"; Kind = "str"; Val = '''' }
            @{ Col = "AW"; Kind = "str"; Val = '''per-erik mukka' }
            @{ Col = "AX"; Kind = "str"; Val = '''per-erik mukka, Anne Järvinen' }
            @{ Col = "AY"; Kind = "str"; Val = '''' }
        )
    }
    @{
        Row = 9
        Cells = @(
            @{ Col = "A"; Kind = "num"; Val = 112501052 }
            @{ Col = "B"; Kind = "num"; Val = 89571 }
            @{ Col = "C"; Kind = "str"; Val = '''Ovaliderad' }
            @{ Col = "D"; Kind = "str"; Val = '''NT' }
            @{ Col = "E"; Kind = "num"; Val = 5432 }
            @{ Col = "F"; Kind = "str"; Val = '''Granticka' }
            @{ Col = "G"; Kind = "str"; Val = '''Porodaedalea chrysoloma' }
            @{ Col = "H"; Kind = "str"; Val = '''(Fr.) Fiasson & Niemelä' }
            @{ Col = "I"; Kind = "str"; Val = '''' }
            @{ Col = "J"; Kind = "str"; Val = '''' }
            @{ Col = "K"; Kind = "str"; Val = '''' }
            @{ Col = "N"; Kind = "str"; Val = '''' }
            @{ Col = "P"; Kind = "str"; Val = '''Erkinjänkkä (Erkinjänkkä), T lm' }
            @{ Col = "Q"; Kind = "num"; Val = 749927 }
            @{ Col = "R"; Kind = "num"; Val = 7535984 }
            @{ Col = "S"; Kind = "num"; Val = 10 }
            @{ Col = "T"; Kind = "str"; Val = '''Norrbotten' }
            @{ Col = "U"; Kind = "str"; Val = '''Kiruna' }
            @{ Col = "V"; Kind = "str"; Val = '''Torne lappmark' }
            @{ Col = "W"; Kind = "str"; Val = '''Jukkasjärvi' }
            @{ Col = "Y"; Kind = "str"; Val = '''2023-10-03' }
            @{ Col = "AA"; Kind = "str"; Val = '''2023-10-03' }
            @{ Col = "AD"; Kind = "bool"; Val = $false }
            @{ Col = "AE"; Kind = "bool"; Val = $false }
            @{ Col = "AF"; Kind = "str"; Val = '''' }
            @{ Col = "AG"; Kind = "bool"; Val = $false }
            @{ Col = "AT"; Kind = "str"; Val = '''' }
            @{ Col = "AW"; Kind = "str"; Val = '''per-erik mukka' }
            @{ Col = "AX"; Kind = "str"; Val = '''per-erik mukka, Anne Järvinen' }
            @{ Col = "AY"; Kind = "str"; Val = '''' }
        )
    }
    @{
        Row = 10
        Cells = @(
            @{ Col = "A"; Kind = "num"; Val = 112501187 }
            @{ Col = "B"; Kind = "num"; Val = 89903 }
            @{ Col = "C"; Kind = "str"; Val = '''Ovaliderad' }
            @{ Col = "D"; Kind = "str"; Val = '''VU' }
            @{ Col = "E"; Kind = "num"; Val = 1506 }
            @{ Col = "F"; Kind = "str"; Val = '''Ostticka' }
            @{ Col = "G"; Kind = "str"; Val = '''Skeletocutis odora' }
            @{ Col = "H"; Kind = "str"; Val = '''(Sacc.) Ginns' }
            @{ Col = "I"; Kind = "str"; Val = '''' }
            @{ Col = "J"; Kind = "str"; Val = '''' }
            @{ Col = "K"; Kind = "str"; Val = '''' }
            @{ Col = "N"; Kind = "str"; Val = '''' }
            @{ Col = "P"; Kind = "str"; Val = '''Erkinjänkkä (Erkinjänkkä), T lm' }
            @{ Col = "Q"; Kind = "num"; Val = 749922 }
            @{ Col = "R"; Kind = "num"; Val = 7535992 }
            @{ Col = "S"; Kind = "num"; Val = 10 }
            @{ Col = "T"; Kind = "str"; Val = '''Norrbotten' }
            @{ Col = "U"; Kind = "str"; Val = '''Kiruna' }
            @{ Col = "V"; Kind = "str"; Val = '''Torne lappmark' }
            @{ Col = "W"; Kind = "str"; Val = '''Jukkasjärvi' }
            @{ Col = "Y"; Kind = "str"; Val = '''2023-10-03' }
            @{ Col = "AA"; Kind = "str"; Val = '''2023-10-03' }
            @{ Col = "AD"; Kind = "bool"; Val = $false }
            @{ Col = "AE"; Kind = "bool"; Val = $false }
            @{ Col = "AF"; Kind = "str"; Val = '''' }
            @{ Col = "AG"; Kind = "bool"; Val = $false }
            @{ Col = "AT"; Kind = "str"; Val = '''' }
            @{ Col = "AW"; Kind = "str"; Val = '''per-erik mukka' }
            @{ Col = "AX"; Kind = "str"; Val = '''per-erik mukka, Anne Järvinen' }
            @{ Col = "AY"; Kind = "str"; Val = '''' }
        )
    }
    @{
        Row = 11
        Cells = @(
            @{ Col = "A"; Kind = "num"; Val = 112501192 }
            @{ Col = "B"; Kind = "num"; Val = 89571 }
            @{ Col = "C"; Kind = "str"; Val = '''Ovaliderad' }
            @{ Col = "D"; Kind = "str"; Val = '''NT' }
            @{ Col = "E"; Kind = "num"; Val = 5432 }
            @{ Col = "F"; Kind = "str"; Val = '''Granticka' }
            @{ Col = "G"; Kind = "str"; Val = '''Porodaedalea chrysoloma' }
            @{ Col = "H"; Kind = "str"; Val = '''(Fr.) Fiasson & Niemelä' }
            @{ Col = "I"; Kind = "str"; Val = '''' }
            @{ Col = "J"; Kind = "str"; Val = '''' }
            @{ Col = "K"; Kind = "str"; Val = '''' }
            @{ Col = "N"; Kind = "str"; Val = '''' }
            @{ Col = "P"; Kind = "str"; Val = '''Erkinjänkkä (Erkinjänkkä), T lm' }
            @{ Col = "Q"; Kind = "num"; Val = 749922 }
            @{ Col = "R"; Kind = "num"; Val = 7535992 }
            @{ Col = "S"; Kind = "num"; Val = 10 }
            @{ Col = "T"; Kind = "str"; Val = '''Norrbotten' }
            @{ Col = "U"; Kind = "str"; Val = '''Kiruna' }
            @{ Col = "V"; Kind = "str"; Val = '''Torne lappmark' }
            @{ Col = "W"; Kind = "str"; Val = '''Jukkasjärvi' }
            @{ Col = "Y"; Kind = "str"; Val = '''2023-10-03' }
            @{ Col = "AA"; Kind = "str"; Val = '''2023-10-03' }
            @{ Col = "AD"; Kind = "bool"; Val = $false }
            @{ Col = "AE"; Kind = "bool"; Val = $false }
            @{ Col = "AF"; Kind = "str"; Val = '''' }
            @{ Col = "AG"; Kind = "bool"; Val = $false }
            @{ Col = "AT"; Kind = "str"; Val = '''' }
            @{ Col = "AW"; Kind = "str"; Val = '''per-erik mukka' }
            @{ Col = "AX"; Kind = "str"; Val = '''per-erik mukka, Anne Järvinen' }
            @{ Col = "AY"; Kind = "str"; Val = '''' }
        )
    }
    @{
        Row = 12
        Cells = @(
            @{ Col = "A"; Kind = "num"; Val = 112501206 }
            @{ Col = "B"; Kind = "num"; Val = 89553 }
            @{ Col = "C"; Kind = "str"; Val = '''Ovaliderad' }
            @{ Col = "D"; Kind = "str"; Val = '''NT' }
            @{ Col = "E"; Kind = "num"; Val = 1202 }
            @{ Col = "F"; Kind = "str"; Val = '''Ullticka' }
            @{ Col = "G"; Kind = "str"; Val = '''Phellinidium ferrugineofuscum' }
            @{ Col = "H"; Kind = "str"; Val = '''(P.Karst.) Fiasson & Niemelä' }
            @{ Col = "I"; Kind = "str"; Val = '''' }
            @{ Col = "J"; Kind = "str"; Val = '''' }
            @{ Col = "K"; Kind = "str"; Val = '''' }
            @{ Col = "N"; Kind = "str"; Val = '''' }
            @{ Col = "P"; Kind = "str"; Val = '''Erkinjänkkä (Erkinjänkkä), T lm' }
            @{ Col = "Q"; Kind = "num"; Val = 749922 }
            @{ Col = "R"; Kind = "num"; Val = 7535992 }
            @{ Col = "S"; Kind = "num"; Val = 10 }
            @{ Col = "T"; Kind = "str"; Val = '''Norrbotten' }
            @{ Col = "U"; Kind = "str"; Val = '''Kiruna' }
            @{ Col = "V"; Kind = "str"; Val = '''Torne lappmark' }
            @{ Col = "W"; Kind = "str"; Val = '''Jukkasjärvi' }
            @{ Col = "Y"; Kind = "str"; Val = '''2023-10-03' }
            @{ Col = "AA"; Kind = "str"; Val = '''2023-10-03' }
            @{ Col = "AD"; Kind = "bool"; Val = $false }
            @{ Col = "AE"; Kind = "bool"; Val = $false }
            @{ Col = "AF"; Kind = "str"; Val = '''' }
            @{ Col = "AG"; Kind = "bool"; Val = $false }
            @{ Col = "AT"; Kind = "str"; Val = '''' }
            @{ Col = "AW"; Kind = "str"; Val = '''per-erik mukka' }
            @{ Col = "AX"; Kind = "str"; Val = '''per-erik mukka, Anne Järvinen' }
            @{ Col = "AY"; Kind = "str"; Val = '''' }
        )
    }
)

foreach ($row in $newRows) {
    foreach ($cell in $row.Cells) {
        $target = $ws.Range("$($cell.Col)$($row.Row)")
        $target.Value = $cell.Val
        # Reset the style after assignment so no implicit number/quote-prefix
        # formatting (e.g. from the date heuristic or the apostrophe prefix)
        # lingers on the cell - the source file uses plain default formatting.
        $target.Style = "Normal"
    }
}
